$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 32,2
$data[0,0] = -0.24293860647387078
$data[0,1] = 0.24227317640985291
$data[1,0] = -0.17643889223892728
$data[1,1] = 0.17461860144539543
$data[2,0] = -0.12489862475423408
$data[2,1] = 0.12430698788621442
$data[3,0] = -0.11630698792661498
$data[3,1] = 0.11576870106182469
$data[4,0] = -0.11276870108592529
$data[4,1] = 0.11093069941905132
$data[5,0] = -0.052743791262694018
$data[5,1] = 0.052125655486879197
$data[6,0] = -0.042125655545563578
$data[6,1] = 0.041969687277267909
$data[7,0] = -0.031969687337991992
$data[7,1] = 0.031676491633244375
$data[8,0] = -0.029676491665536098
$data[8,1] = 0.029430142313080054
$data[9,0] = -0.027430142347501629
$data[9,1] = 0.027413225340627179
$data[10,0] = -0.024413225379130488
$data[10,1] = 0.024385288048343945
$data[11,0] = -0.020885288089329546
$data[11,1] = 0.020680462168908065
$data[12,0] = -0.017180462212100345
$data[12,1] = 0.01708706853581976
$data[13,0] = -0.009087068596918435
$data[13,1] = 0.0090564476625312196
$data[14,0] = -0.0080564476970916843
$data[14,1] = 0.0080366634581743313
$data[15,0] = -0.0060366634969679644
$data[15,1] = 0.006003779790713093
$data[16,0] = -0.004003779830023646
$data[16,1] = 0.0039999999528044228
$data[17,0] = -0.029719983353512447
$data[17,1] = 0.029679483905905357
$data[18,0] = -0.025679483922766533
$data[18,1] = 0.025413744550709882
$data[19,0] = -0.008017068953023454
$data[19,1] = 0.008005718968771447
$data[20,0] = -0.0040057189881821387
$data[20,1] = 0.0039999999804347652
$data[21,0] = -0.045719976713353816
$data[21,1] = 0.045504178974223208
$data[22,0] = -0.040504179001483287
$data[22,1] = 0.040100089208925382
$data[23,0] = -0.020100089296197332
$data[23,1] = 0.019999999911560984
$data[24,0] = -0.056186908183653017
$data[24,1] = 0.056151268951049715
$data[25,0] = -0.053651268980898337
$data[25,1] = 0.053608216882031101
$data[26,0] = -0.051108216912754134
$data[26,1] = 0.050869249410070871
$data[27,0] = -0.048869249442774709
$data[27,1] = 0.04871744571833414
$data[28,0] = -0.041717445773058692
$data[28,1] = 0.041682425136279733
$data[29,0] = 0.018317574603930797
$data[29,1] = -0.01834130520180155
$data[30,0] = -0.014023717826574966
$data[30,1] = 0.014001351670700046
$data[31,0] = -0.0040013517374859475
$data[31,1] = 0.0039999999561963762

$ws.Range("A1:B32").Value = $data

# Column B width changes from 15.42578125 to 14.7109375 characters.
# The COM ColumnWidth setter here quantizes to the nearest 1/6-character
# (pixel) increment, so 13.85 is the input that lands on the closest
# representable stored width (14.666666666666666 -> rounds visually to 14.7).
$ws.Columns.Item(2).ColumnWidth = 13.85
